{"js": "// Ordered [oldText, newText] pairs for every paragraph in the document body,\n// in document order: the title line followed by each table cell (row-major).\nconst REPLACEMENTS = [\n  [\"2024-05-30 Thursday\", \"2024-05-31 Friday\"],\n  [\"62-25=37\", \"37+36=73\"],\n  [\"57+36=93\", \"6+7=13\"],\n  [\"7+76=83\", \"44-38=6\"],\n  [\"81-28=53\", \"62+19=81\"],\n  [\"82-53=29\", \"57+39=96\"],\n  [\"73-67=6\", \"18+67=85\"],\n  [\"47+28=75\", \"70-38=32\"],\n  [\"54+8=62\", \"6+78=84\"],\n  [\"38+35=73\", \"74-26=48\"],\n  [\"59+33=92\", \"54-47=7\"],\n  [\"47+19=66\", \"39+27=66\"],\n  [\"49+32=81\", \"83-29=54\"],\n  [\"8+7=15\", \"54+7=61\"],\n  [\"22+39=61\", \"60-9=51\"],\n  [\"55-27=28\", \"84-47=37\"],\n  [\"77+19=96\", \"9+64=73\"],\n  [\"33+19=52\", \"90-6=84\"],\n  [\"31-4=27\", \"74-58=16\"],\n  [\"74-57=17\", \"69+22=91\"],\n  [\"60-27=33\", \"4+27=31\"],\n  [\"98-59=39\", \"33-25=8\"],\n  [\"10-9=1\", \"80-58=22\"],\n  [\"66-29=37\", \"58+37=95\"],\n  [\"5+59=64\", \"4+78=82\"],\n  [\"82-63=19\", \"81-56=25\"],\n  [\"90-19=71\", \"62-53=9\"],\n  [\"19+45=64\", \"93-44=49\"],\n  [\"60-53=7\", \"63+18=81\"],\n  [\"95-29=66\", \"59+24=83\"],\n  [\"43-19=24\", \"64-6=58\"],\n  [\"44-28=16\", \"28+36=64\"],\n  [\"19+68=87\", \"71-46=25\"],\n  [\"17+24=41\", \"76+9=85\"],\n  [\"85-49=36\", \"20-7=13\"],\n  [\"38+5=43\", \"93-25=68\"],\n  [\"44-26=18\", \"66-9=57\"],\n  [\"15+27=42\", \"17+38=55\"],\n  [\"43+28=71\", \"62-14=48\"],\n  [\"94-15=79\", \"95-28=67\"],\n  [\"41-24=17\", \"62-26=36\"],\n  [\"28-19=9\", \"73-44=29\"],\n  [\"45-39=6\", \"48+16=64\"],\n  [\"19+56=75\", \"60-58=2\"],\n  [\"62-9=53\", \"46-8=38\"],\n  [\"8+33=41\", \"49+3=52\"],\n  [\"87-9=78\", \"63-36=27\"],\n  [\"77-18=59\", \"98-19=79\"],\n  [\"20-16=4\", \"26+48=74\"],\n  [\"38+59=97\", \"29+27=56\"],\n  [\"68+13=81\", \"36+59=95\"],\n  [\"91-5=86\", \"91-87=4\"],\n  [\"74+17=91\", \"38+29=67\"],\n  [\"44+49=93\", \"29+38=67\"],\n  [\"60-34=26\", \"76+16=92\"],\n  [\"63-44=19\", \"68+16=84\"],\n  [\"55-27=28\", \"52-38=14\"],\n  [\"50-8=42\", \"7+49=56\"],\n  [\"77+19=96\", \"62-14=48\"],\n  [\"38+25=63\", \"66+17=83\"],\n  [\"9+82=91\", \"30-23=7\"],\n  [\"14-8=6\", \"82-79=3\"],\n  [\"94-75=19\", \"17+75=92\"],\n  [\"94-79=15\", \"81-24=57\"],\n  [\"48+44=92\", \"92-89=3\"],\n  [\"80-56=24\", \"25+46=71\"],\n  [\"95-27=68\", \"38+28=66\"],\n  [\"27+39=66\", \"39+45=84\"],\n  [\"64-58=6\", \"73-57=16\"],\n  [\"19+12=31\", \"9+2=11\"],\n  [\"80-79=1\", \"28+24=52\"],\n  [\"50-6=44\", \"46-17=29\"],\n  [\"27+47=74\", \"37+9=46\"],\n  [\"92-76=16\", \"13+28=41\"],\n  [\"16+29=45\", \"74-46=28\"],\n  [\"33-16=17\", \"58+33=91\"],\n  [\"17+9=26\", \"80-74=6\"],\n  [\"72-34=38\", \"47+26=73\"],\n  [\"31-3=28\", \"66+25=91\"],\n  [\"41-33=8\", \"28+29=57\"],\n  [\"86-18=68\", \"67+27=94\"],\n  [\"73-37=36\", \"67-29=38\"],\n  [\"59+9=68\", \"6+27=33\"],\n  [\"95-69=26\", \"35+19=54\"],\n  [\"25+48=73\", \"15+78=93\"],\n  [\"96-28=68\", \"56-39=17\"],\n  [\"85-36=49\", \"64-36=28\"],\n  [\"62-49=13\", \"28+17=45\"],\n  [\"48+7=55\", \"93-24=69\"],\n  [\"91-19=72\", \"61-27=34\"],\n  [\"35+26=61\", \"9+3=12\"],\n  [\"7+9=16\", \"8+76=84\"],\n  [\"26+66=92\", \"17+68=85\"],\n  [\"31-29=2\", \"5+89=94\"],\n  [\"45+8=53\", \"60-41=19\"],\n  [\"61-34=27\", \"28+17=45\"],\n  [\"62-27=35\", \"60-48=12\"],\n  [\"18+55=73\", \"4+18=22\"],\n  [\"49+42=91\", \"53-5=48\"],\n  [\"25+37=62\", \"92-6=86\"],\n  [\"21-19=2\", \"47+17=64\"]\n];\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst items = paragraphs.items;\nif (items.length !== REPLACEMENTS.length) {\n  throw new Error(\n    \"Unexpected paragraph count: \" + items.length + \" vs expected \" + REPLACEMENTS.length\n  );\n}\n\nfor (let i = 0; i < items.length; i++) {\n  const [oldText, newText] = REPLACEMENTS[i];\n  const current = items[i].text.replace(/[\\r\\x07]+$/g, \"\");\n  if (current !== oldText) {\n    throw new Error(\n      \"Mismatch at paragraph \" + i + \": expected '\" + oldText + \"' but found '\" + current + \"'\"\n    );\n  }\n  if (oldText !== newText) {\n    items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "$oldTexts = @(\n    \"2024-05-30 Thursday\",\n    \"62-25=37\",\n    \"57+36=93\",\n    \"7+76=83\",\n    \"81-28=53\",\n    \"82-53=29\",\n    \"73-67=6\",\n    \"47+28=75\",\n    \"54+8=62\",\n    \"38+35=73\",\n    \"59+33=92\",\n    \"47+19=66\",\n    \"49+32=81\",\n    \"8+7=15\",\n    \"22+39=61\",\n    \"55-27=28\",\n    \"77+19=96\",\n    \"33+19=52\",\n    \"31-4=27\",\n    \"74-57=17\",\n    \"60-27=33\",\n    \"98-59=39\",\n    \"10-9=1\",\n    \"66-29=37\",\n    \"5+59=64\",\n    \"82-63=19\",\n    \"90-19=71\",\n    \"19+45=64\",\n    \"60-53=7\",\n    \"95-29=66\",\n    \"43-19=24\",\n    \"44-28=16\",\n    \"19+68=87\",\n    \"17+24=41\",\n    \"85-49=36\",\n    \"38+5=43\",\n    \"44-26=18\",\n    \"15+27=42\",\n    \"43+28=71\",\n    \"94-15=79\",\n    \"41-24=17\",\n    \"28-19=9\",\n    \"45-39=6\",\n    \"19+56=75\",\n    \"62-9=53\",\n    \"8+33=41\",\n    \"87-9=78\",\n    \"77-18=59\",\n    \"20-16=4\",\n    \"38+59=97\",\n    \"68+13=81\",\n    \"91-5=86\",\n    \"74+17=91\",\n    \"44+49=93\",\n    \"60-34=26\",\n    \"63-44=19\",\n    \"55-27=28\",\n    \"50-8=42\",\n    \"77+19=96\",\n    \"38+25=63\",\n    \"9+82=91\",\n    \"14-8=6\",\n    \"94-75=19\",\n    \"94-79=15\",\n    \"48+44=92\",\n    \"80-56=24\",\n    \"95-27=68\",\n    \"27+39=66\",\n    \"64-58=6\",\n    \"19+12=31\",\n    \"80-79=1\",\n    \"50-6=44\",\n    \"27+47=74\",\n    \"92-76=16\",\n    \"16+29=45\",\n    \"33-16=17\",\n    \"17+9=26\",\n    \"72-34=38\",\n    \"31-3=28\",\n    \"41-33=8\",\n    \"86-18=68\",\n    \"73-37=36\",\n    \"59+9=68\",\n    \"95-69=26\",\n    \"25+48=73\",\n    \"96-28=68\",\n    \"85-36=49\",\n    \"62-49=13\",\n    \"48+7=55\",\n    \"91-19=72\",\n    \"35+26=61\",\n    \"7+9=16\",\n    \"26+66=92\",\n    \"31-29=2\",\n    \"45+8=53\",\n    \"61-34=27\",\n    \"62-27=35\",\n    \"18+55=73\",\n    \"49+42=91\",\n    \"25+37=62\",\n    \"21-19=2\"\n)\n$newTexts = @(\n    \"2024-05-31 Friday\",\n    \"37+36=73\",\n    \"6+7=13\",\n    \"44-38=6\",\n    \"62+19=81\",\n    \"57+39=96\",\n    \"18+67=85\",\n    \"70-38=32\",\n    \"6+78=84\",\n    \"74-26=48\",\n    \"54-47=7\",\n    \"39+27=66\",\n    \"83-29=54\",\n    \"54+7=61\",\n    \"60-9=51\",\n    \"84-47=37\",\n    \"9+64=73\",\n    \"90-6=84\",\n    \"74-58=16\",\n    \"69+22=91\",\n    \"4+27=31\",\n    \"33-25=8\",\n    \"80-58=22\",\n    \"58+37=95\",\n    \"4+78=82\",\n    \"81-56=25\",\n    \"62-53=9\",\n    \"93-44=49\",\n    \"63+18=81\",\n    \"59+24=83\",\n    \"64-6=58\",\n    \"28+36=64\",\n    \"71-46=25\",\n    \"76+9=85\",\n    \"20-7=13\",\n    \"93-25=68\",\n    \"66-9=57\",\n    \"17+38=55\",\n    \"62-14=48\",\n    \"95-28=67\",\n    \"62-26=36\",\n    \"73-44=29\",\n    \"48+16=64\",\n    \"60-58=2\",\n    \"46-8=38\",\n    \"49+3=52\",\n    \"63-36=27\",\n    \"98-19=79\",\n    \"26+48=74\",\n    \"29+27=56\",\n    \"36+59=95\",\n    \"91-87=4\",\n    \"38+29=67\",\n    \"29+38=67\",\n    \"76+16=92\",\n    \"68+16=84\",\n    \"52-38=14\",\n    \"7+49=56\",\n    \"62-14=48\",\n    \"66+17=83\",\n    \"30-23=7\",\n    \"82-79=3\",\n    \"17+75=92\",\n    \"81-24=57\",\n    \"92-89=3\",\n    \"25+46=71\",\n    \"38+28=66\",\n    \"39+45=84\",\n    \"73-57=16\",\n    \"9+2=11\",\n    \"28+24=52\",\n    \"46-17=29\",\n    \"37+9=46\",\n    \"13+28=41\",\n    \"74-46=28\",\n    \"58+33=91\",\n    \"80-74=6\",\n    \"47+26=73\",\n    \"66+25=91\",\n    \"28+29=57\",\n    \"67+27=94\",\n    \"67-29=38\",\n    \"6+27=33\",\n    \"35+19=54\",\n    \"15+78=93\",\n    \"56-39=17\",\n    \"64-36=28\",\n    \"28+17=45\",\n    \"93-24=69\",\n    \"61-27=34\",\n    \"9+3=12\",\n    \"8+76=84\",\n    \"17+68=85\",\n    \"5+89=94\",\n    \"60-41=19\",\n    \"28+17=45\",\n    \"60-48=12\",\n    \"4+18=22\",\n    \"53-5=48\",\n    \"92-6=86\",\n    \"47+17=64\"\n)\n\n$d = $word.ActiveDocument\n$count = $d.Paragraphs.Count\n\n$idx = 0\nfor ($i = 1; $i -le $count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    $r = $p.Range\n    $t = $r.Text\n    $clean = $t -replace \"[\\x07\\r\\n]\", \"\"\n    if ($clean -eq \"\") {\n        continue\n    }\n    if ($idx -ge $oldTexts.Length) {\n        break\n    }\n    if ($clean -ne $oldTexts[$idx]) {\n        throw \"Mismatch at paragraph $i (item $idx): expected '$($oldTexts[$idx])' got '$clean'\"\n    }\n    if ($oldTexts[$idx] -ne $newTexts[$idx]) {\n        $r.Text = $newTexts[$idx]\n    }\n    $idx = $idx + 1\n}\n\nif ($idx -ne $oldTexts.Length) {\n    throw \"Only processed $idx of $($oldTexts.Length) replacements\"\n}\n"}
